$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update stats for 2025-10 (row 23)
$ws.Range("B23").Value = 6331
$ws.Range("C23").Value = 1000
$ws.Range("D23").Value = 5905747
$ws.Range("E23").Value = 932.8300426472911
$ws.Range("F23").Value = 8.63074811256006
$ws.Range("G23").Value = 4.058272632674287
$ws.Range("H23").Value = 26.54504165815819
